# Update "想去人数" (want-to-go count) values in column F across the
# four worksheets, as captured by the source diff.

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1179
$ws1.Range("F5").Value = 320
$ws1.Range("F10").Value = 299
$ws1.Range("F13").Value = 119
$ws1.Range("F19").Value = 263

# 演出 (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F10").Value = 68
$ws2.Range("F13").Value = 216
$ws2.Range("F23").Value = 40

# 本地生活 (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 801
$ws3.Range("F4").Value = 2041

# 全部类型 (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 801
$ws4.Range("F4").Value = 2041
$ws4.Range("F13").Value = 1179
$ws4.Range("F14").Value = 320
$ws4.Range("F21").Value = 68
$ws4.Range("F26").Value = 299
$ws4.Range("F27").Value = 216
$ws4.Range("F30").Value = 119
$ws4.Range("F45").Value = 40
$ws4.Range("F46").Value = 263
